$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 223.23077
$ws.Range("J2").Value = 273.33334
$ws.Range("L2").Value = 273.33334
$ws.Range("N2").Value = -499.33334

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 686.86957
$ws.Range("I28").Value = 661.2941
$ws.Range("J28").Value = 759.3333
$ws.Range("K28").Value = 661.2941
$ws.Range("L28").Value = 759.3333
$ws.Range("M28").Value = -176.2941
$ws.Range("N28").Value = -1729.3333

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3582.5518
$ws.Range("I62").Value = 2688.611
$ws.Range("J62").Value = 5045.364
$ws.Range("K62").Value = 2688.611
$ws.Range("L62").Value = 5045.364
$ws.Range("M62").Value = -2064.611
$ws.Range("N62").Value = -6293.364

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3582.5518
$ws.Range("I65").Value = 2688.611
$ws.Range("J65").Value = 5045.364
$ws.Range("K65").Value = 13443.055
$ws.Range("L65").Value = 25226.82
$ws.Range("M65").Value = -10323.055
$ws.Range("N65").Value = -31466.82

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3761.389
$ws.Range("J69").Value = 4696.875
$ws.Range("L69").Value = 14090.625
$ws.Range("N69").Value = -15838.625

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 3761.389
$ws.Range("J72").Value = 4696.875
$ws.Range("L72").Value = 42271.875
$ws.Range("N72").Value = -51007.875

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3424.878
$ws.Range("I76").Value = 2981.818
$ws.Range("K76").Value = 2981.818
$ws.Range("M76").Value = -2666.818

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3424.878
$ws.Range("I79").Value = 2981.818
$ws.Range("K79").Value = 2981.818
$ws.Range("M79").Value = -1889.818

# ALC row 115
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 793.7646999999999
$ws.Range("I115").Value = 224.25
$ws.Range("K115").Value = 672.75
$ws.Range("M115").Value = 894.25

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1494.6842
$ws.Range("J121").Value = 1680.2667
$ws.Range("L121").Value = 5040.800099999999
$ws.Range("N121").Value = -8534.8001

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1018.9211
$ws.Range("J129").Value = 1054.8429
$ws.Range("L129").Value = 3164.5287
$ws.Range("N129").Value = -13164.5287

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3775.31
$ws.Range("I138").Value = 2658.2058
$ws.Range("J138").Value = 4350.788
$ws.Range("K138").Value = 7974.617400000001
$ws.Range("L138").Value = 13052.364
$ws.Range("M138").Value = -2834.617400000001
$ws.Range("N138").Value = -23332.364

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11166.395
$ws.Range("I32").Value = 8826.099
$ws.Range("J32").Value = 44398.6
$ws.Range("K32").Value = 8826.099
$ws.Range("L32").Value = 44398.6
$ws.Range("M32").Value = -8539.099
$ws.Range("N32").Value = -44972.6

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4000
$ws.Range("J63").Value = 4000
$ws.Range("L63").Value = 4000
$ws.Range("N63").Value = -5372

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4000
$ws.Range("J66").Value = 4000
$ws.Range("L66").Value = 20000
$ws.Range("N66").Value = -26864

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2759.4707
$ws.Range("I122").Value = 1731.1
$ws.Range("J122").Value = 4228.5713
$ws.Range("K122").Value = 5193.299999999999
$ws.Range("L122").Value = 12685.7139
$ws.Range("M122").Value = -2743.299999999999
$ws.Range("N122").Value = -17585.7139

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1774.1875
$ws.Range("I132").Value = 1417.3273
$ws.Range("J132").Value = 2559.28
$ws.Range("K132").Value = 4251.9819
$ws.Range("L132").Value = 7677.84
$ws.Range("M132").Value = -1721.9819
$ws.Range("N132").Value = -12737.84

# BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 31433.334
$ws.Range("I35").Value = 20000
$ws.Range("J35").Value = 33720
$ws.Range("K35").Value = 20000
$ws.Range("L35").Value = 33720
$ws.Range("M35").Value = -19690
$ws.Range("N35").Value = -34340

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 8900
$ws.Range("I82").Value = 8900
$ws.Range("K82").Value = 8900
$ws.Range("M82").Value = -8517

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 8900
$ws.Range("I85").Value = 8900
$ws.Range("K85").Value = 8900
$ws.Range("M85").Value = -7574

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1728.2222
$ws.Range("I134").Value = 1261.8438
$ws.Range("J134").Value = 2876.2307
$ws.Range("K134").Value = 3785.5314
$ws.Range("L134").Value = 8628.6921
$ws.Range("M134").Value = -1250.5314
$ws.Range("N134").Value = -13698.6921

# CRP row 57
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 50000
$ws.Range("J57").Value = 50000
$ws.Range("L57").Value = 50000
$ws.Range("N57").Value = -51120

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4438.4375
$ws.Range("I62").Value = 4101.8
$ws.Range("J62").Value = 4591.4546
$ws.Range("K62").Value = 4101.8
$ws.Range("L62").Value = 4591.4546
$ws.Range("M62").Value = -3477.8
$ws.Range("N62").Value = -5839.4546

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4438.4375
$ws.Range("I65").Value = 4101.8
$ws.Range("J65").Value = 4591.4546
$ws.Range("K65").Value = 20509
$ws.Range("L65").Value = 22957.273
$ws.Range("M65").Value = -17389
$ws.Range("N65").Value = -29197.273

# CUL row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 226.15384
$ws.Range("J40").Value = 304.44446
$ws.Range("L40").Value = 1217.77784
$ws.Range("N40").Value = -1355.77784

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1110.89
$ws.Range("I131").Value = 435.83334
$ws.Range("J131").Value = 1202.9432
$ws.Range("K131").Value = 1307.50002
$ws.Range("L131").Value = 3608.8296
$ws.Range("M131").Value = 3732.49998
$ws.Range("N131").Value = -13688.8296

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2563.7856
$ws.Range("I132").Value = 1999.3334
$ws.Range("K132").Value = 5998.0002
$ws.Range("M132").Value = -3468.0002

# GSM row 133
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 92780
$ws.Range("J133").Value = 92780
$ws.Range("L133").Value = 92780
$ws.Range("N133").Value = -102900

# LTW row 39
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 260.92307
$ws.Range("I55").Value = 160.2
$ws.Range("J55").Value = 323.875
$ws.Range("K55").Value = 160.2
$ws.Range("L55").Value = 323.875
$ws.Range("M55").Value = 12.80000000000001
$ws.Range("N55").Value = -669.875

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1548.625
$ws.Range("I93").Value = 1163.3334
$ws.Range("J93").Value = 1779.8
$ws.Range("K93").Value = 1163.3334
$ws.Range("L93").Value = 1779.8
$ws.Range("M93").Value = 84.66660000000002
$ws.Range("N93").Value = -4275.8

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2485.8948
$ws.Range("I136").Value = 2485.8948
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7457.6844
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4907.6844
$ws.Range("N136").ClearContents()

# WVR row 44
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 49999.5
$ws.Range("J44").Value = 49999.5
$ws.Range("L44").Value = 49999.5
$ws.Range("N44").Value = -51107.5

# WVR row 58
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 6000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 6000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -6616

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3552.5667
$ws.Range("I136").Value = 671.48486
$ws.Range("J136").Value = 7073.8887
$ws.Range("K136").Value = 2014.45458
$ws.Range("L136").Value = 21221.6661
$ws.Range("M136").Value = 535.5454199999999
$ws.Range("N136").Value = -26321.6661
